# Applies classroom re-allocation changes described in the commit diff.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Section_A sheet: update elective room labels
# ---------------------------------------------------------------------
$wsA = $wb.Worksheets.Item("Section_A")
$wsA.Range("B2").Value = "ELECTIVE_B6 [C102]"
$wsA.Range("C2").Value = "ELECTIVE_B7 [C403]"
$wsA.Range("D5").Value = "ELECTIVE_B6 [C102]"
$wsA.Range("E5").Value = "ELECTIVE_B7 [C403]"
$wsA.Range("C6").Value = "ELECTIVE_B6 (Tutorial) [C303]"
$wsA.Range("D6").Value = "ELECTIVE_B7 (Tutorial) [C205]"

# ---------------------------------------------------------------------
# Section_B sheet: update elective room labels
# ---------------------------------------------------------------------
$wsB = $wb.Worksheets.Item("Section_B")
$wsB.Range("B2").Value = "ELECTIVE_B6 [C405]"
$wsB.Range("C2").Value = "ELECTIVE_B7 [C004]"
$wsB.Range("D5").Value = "ELECTIVE_B6 [C405]"
$wsB.Range("E5").Value = "ELECTIVE_B7 [C004]"
$wsB.Range("C6").Value = "ELECTIVE_B6 (Tutorial) [C305]"
$wsB.Range("D6").Value = "ELECTIVE_B7 (Tutorial) [C305]"

# ---------------------------------------------------------------------
# Classroom_Utilization sheet: update weekly hours / daily avg / utilization
# for the rooms whose bookings moved.
# ---------------------------------------------------------------------
$wsU = $wb.Worksheets.Item("Classroom_Utilization")

# C002 (row 3): now unused
$wsU.Range("D3").Value = 0
$wsU.Range("E3").Value = 0
$wsU.Range("G3").Value = 0

# C004 (row 5): now used
$wsU.Range("D5").Value = 3
$wsU.Range("E5").Value = 0.6
$wsU.Range("G5").Value = 7.5

# C102 (row 7): now used
$wsU.Range("D7").Value = 3
$wsU.Range("E7").Value = 0.6
$wsU.Range("G7").Value = 7.5

# C104 (row 9): now unused
$wsU.Range("D9").Value = 0
$wsU.Range("E9").Value = 0
$wsU.Range("G9").Value = 0

# C203 (row 15): now unused
$wsU.Range("D15").Value = 0
$wsU.Range("E15").Value = 0
$wsU.Range("G15").Value = 0

# C204 (row 16): now unused
$wsU.Range("D16").Value = 0
$wsU.Range("E16").Value = 0
$wsU.Range("G16").Value = 0

# C205 (row 17): now used
$wsU.Range("D17").Value = 1
$wsU.Range("E17").Value = 0.2
$wsU.Range("G17").Value = 2.5

# C305 (row 25): now used more
$wsU.Range("D25").Value = 2
$wsU.Range("E25").Value = 0.4
$wsU.Range("G25").Value = 5

# C402 (row 30): now unused
$wsU.Range("D30").Value = 0
$wsU.Range("E30").Value = 0
$wsU.Range("G30").Value = 0

# C403 (row 31): now used
$wsU.Range("D31").Value = 3
$wsU.Range("E31").Value = 0.6
$wsU.Range("G31").Value = 7.5

# C405 (row 33): usage decreased
$wsU.Range("D33").Value = 3
$wsU.Range("E33").Value = 0.6
$wsU.Range("G33").Value = 7.5

# ---------------------------------------------------------------------
# Classroom_Allocation sheet: update room assignments (and derived
# room type / capacity columns) for each scheduled session.
# The "Capacity" column (I) is stored as text in the workbook, so force
# a text number format before assigning the numeric-looking strings —
# otherwise the value would be auto-coerced back into a number.
# ---------------------------------------------------------------------
$wsC = $wb.Worksheets.Item("Classroom_Allocation")
$capacityCells = @("I2","I3","I5","I7","I8","I9","I10","I11","I13")
foreach ($addr in $capacityCells) {
    $wsC.Range($addr).NumberFormat = "@"
}

# Row 2: Section A, ELECTIVE_B6 Mon 09:00-10:30
$wsC.Range("G2").Value = "C102"
$wsC.Range("I2").Value = "96"

# Row 3: Section A, ELECTIVE_B7 Tue 09:00-10:30
$wsC.Range("G3").Value = "C403"
$wsC.Range("I3").Value = "78"

# Row 4: Section A, ELECTIVE_B6 (Tutorial) Tue 14:30-15:30
$wsC.Range("G4").Value = "C303"

# Row 5: Section A, ELECTIVE_B6 Wed 13:00-14:30
$wsC.Range("G5").Value = "C102"
$wsC.Range("I5").Value = "96"

# Row 6: Section A, ELECTIVE_B7 (Tutorial) Wed 14:30-15:30
$wsC.Range("G6").Value = "C205"

# Row 7: Section A, ELECTIVE_B7 Thu 13:00-14:30
$wsC.Range("G7").Value = "C403"
$wsC.Range("I7").Value = "78"

# Row 8: Section B, ELECTIVE_B6 Mon 09:00-10:30
$wsC.Range("G8").Value = "C405"
$wsC.Range("H8").Value = "classroom"
$wsC.Range("I8").Value = "78"

# Row 9: Section B, ELECTIVE_B7 Tue 09:00-10:30
$wsC.Range("G9").Value = "C004"
$wsC.Range("H9").Value = "Auditorium"
$wsC.Range("I9").Value = "240"

# Row 10: Section B, ELECTIVE_B6 (Tutorial) Tue 14:30-15:30
$wsC.Range("G10").Value = "C305"
$wsC.Range("I10").Value = "96"

# Row 11: Section B, ELECTIVE_B6 Wed 13:00-14:30
$wsC.Range("G11").Value = "C405"
$wsC.Range("H11").Value = "classroom"
$wsC.Range("I11").Value = "78"

# Row 12: Section B, ELECTIVE_B7 (Tutorial) Wed 14:30-15:30
$wsC.Range("G12").Value = "C305"

# Row 13: Section B, ELECTIVE_B7 Thu 13:00-14:30
$wsC.Range("G13").Value = "C004"
$wsC.Range("H13").Value = "Auditorium"
$wsC.Range("I13").Value = "240"

Write-Output "Applied classroom re-allocation updates."
